$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 223.33333
$ws.Range("I2").Value = 223.33333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 223.33333
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -110.33333

$ws.Range("H6").Value = 936.9
$ws.Range("I6").Value = 146.35294
$ws.Range("J6").Value = 5416.6665
$ws.Range("K6").Value = 439.05882
$ws.Range("L6").Value = 16249.9995
$ws.Range("M6").Value = -327.05882
$ws.Range("N6").Value = -16473.9995

$ws.Range("H17").Value = 1956
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1956
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5868
$ws.Range("N17").Value = -6204

$ws.Range("H64").Value = 4637.75
$ws.Range("I64").Value = 2501
$ws.Range("J64").Value = 5350
$ws.Range("K64").Value = 2501
$ws.Range("L64").Value = 5350
$ws.Range("M64").Value = -2253
$ws.Range("N64").Value = -5846

$ws.Range("H67").Value = 4637.75
$ws.Range("I67").Value = 2501
$ws.Range("J67").Value = 5350
$ws.Range("K67").Value = 2501
$ws.Range("L67").Value = 5350
$ws.Range("M67").Value = -1643
$ws.Range("N67").Value = -7066

$ws.Range("H74").Value = 3999
$ws.Range("I74").Value = 3999
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3999
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3063

$ws.Range("H77").Value = 3999
$ws.Range("I77").Value = 3999
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 19995
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -15315

$ws.Range("H81").Value = 40000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 40000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996

$ws.Range("H84").Value = 40000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 40000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984

$ws.Range("H137").Value = 1352.6666
$ws.Range("I137").Value = 1216.875
$ws.Range("J137").Value = 1624.25
$ws.Range("K137").Value = 3650.625
$ws.Range("L137").Value = 4872.75
$ws.Range("M137").Value = -1100.625
$ws.Range("N137").Value = -9972.75

$ws.Range("H138").Value = 5143.4883
$ws.Range("I138").Value = 2764.432
$ws.Range("J138").Value = 7635.8335
$ws.Range("K138").Value = 8293.295999999998
$ws.Range("L138").Value = 22907.5005
$ws.Range("M138").Value = -3153.295999999998
$ws.Range("N138").Value = -33187.50049999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 71.28570999999999
$ws.Range("I5").Value = 20
$ws.Range("J5").Value = 139.66667
$ws.Range("K5").Value = 20
$ws.Range("L5").Value = 139.66667
$ws.Range("M5").Value = 92
$ws.Range("N5").Value = -363.66667

$ws.Range("H61").Value = 2558.6365
$ws.Range("I61").Value = 2558.6365
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2558.6365
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2346.6365

$ws.Range("H97").Value = 1163.5555
$ws.Range("I97").Value = 1261.25
$ws.Range("J97").Value = 382
$ws.Range("K97").Value = 1261.25
$ws.Range("L97").Value = 382
$ws.Range("M97").Value = -765.25
$ws.Range("N97").Value = -1374

$ws.Range("H112").Value = 38746.75
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 38746.75
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 38746.75
$ws.Range("N112").Value = -41700.75

$ws.Range("H122").Value = 1706
$ws.Range("I122").Value = 1706
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5118
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2668

$ws.Range("H132").Value = 1016.48
$ws.Range("I132").Value = 975.5
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2926.5
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -396.5
$ws.Range("N132").Value = -11060

$ws.Range("H136").Value = 2558.6365
$ws.Range("I136").Value = 2558.6365
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7675.9095
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -5125.9095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 71.28570999999999
$ws.Range("I4").Value = 20
$ws.Range("J4").Value = 139.66667
$ws.Range("K4").Value = 20
$ws.Range("L4").Value = 139.66667
$ws.Range("M4").Value = 95
$ws.Range("N4").Value = -369.66667

$ws.Range("H20").Value = 4649
$ws.Range("I20").Value = 3897.9092
$ws.Range("J20").Value = 7403
$ws.Range("K20").Value = 3897.9092
$ws.Range("L20").Value = 7403
$ws.Range("M20").Value = -3650.9092
$ws.Range("N20").Value = -7897

$ws.Range("H22").Value = 118051.766
$ws.Range("I22").Value = 462.92856
$ws.Range("J22").Value = 666799.7
$ws.Range("K22").Value = 462.92856
$ws.Range("L22").Value = 666799.7
$ws.Range("M22").Value = -289.92856
$ws.Range("N22").Value = -667145.7

$ws.Range("H86").Value = 4359.091
$ws.Range("I86").Value = 3272.2222
$ws.Range("J86").Value = 9250
$ws.Range("K86").Value = 3272.2222
$ws.Range("L86").Value = 9250
$ws.Range("M86").Value = -2149.2222
$ws.Range("N86").Value = -11496

$ws.Range("H89").Value = 4359.091
$ws.Range("I89").Value = 3272.2222
$ws.Range("J89").Value = 9250
$ws.Range("K89").Value = 16361.111
$ws.Range("L89").Value = 46250
$ws.Range("M89").Value = -10745.111
$ws.Range("N89").Value = -57482

$ws.Range("H107").Value = 563.8182
$ws.Range("I107").Value = 404.875
$ws.Range("J107").Value = 987.6667
$ws.Range("K107").Value = 404.875
$ws.Range("L107").Value = 987.6667
$ws.Range("M107").Value = 1515.125
$ws.Range("N107").Value = -4827.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 60535
$ws.Range("I62").Value = 3249.5
$ws.Range("J62").Value = 83449.2
$ws.Range("K62").Value = 3249.5
$ws.Range("L62").Value = 83449.2
$ws.Range("M62").Value = -2625.5
$ws.Range("N62").Value = -84697.2

$ws.Range("H65").Value = 60535
$ws.Range("I65").Value = 3249.5
$ws.Range("J65").Value = 83449.2
$ws.Range("K65").Value = 16247.5
$ws.Range("L65").Value = 417246
$ws.Range("M65").Value = -13127.5
$ws.Range("N65").Value = -423486

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 88966.664
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 88966.664
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 266899.992
$ws.Range("N37").Value = -267123.992

$ws.Range("H103").Value = 168.6
$ws.Range("I103").Value = 178.75
$ws.Range("J103").Value = 128
$ws.Range("K103").Value = 536.25
$ws.Range("L103").Value = 384
$ws.Range("M103").Value = 342.75
$ws.Range("N103").Value = -2142

$ws.Range("H113").Value = 2634.4
$ws.Range("I113").Value = 2786
$ws.Range("J113").Value = 2533.3333
$ws.Range("K113").Value = 8358
$ws.Range("L113").Value = 7599.999899999999
$ws.Range("M113").Value = -6188
$ws.Range("N113").Value = -11939.9999

$ws.Range("H128").Value = 2196942
$ws.Range("I128").Value = 2196942
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 6590826
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -6585846

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 573.7857
$ws.Range("I2").Value = 48
$ws.Range("J2").Value = 968.125
$ws.Range("K2").Value = 48
$ws.Range("L2").Value = 968.125
$ws.Range("M2").Value = 65
$ws.Range("N2").Value = -1194.125

$ws.Range("H43").Value = 1218.0714
$ws.Range("I43").Value = 587.75
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 587.75
$ws.Range("L43").Value = 5000
$ws.Range("M43").Value = -436.75
$ws.Range("N43").Value = -5302

$ws.Range("H109").Value = 48000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 48000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 48000
$ws.Range("N109").Value = -50080

$ws.Range("H113").Value = 1999
$ws.Range("I113").Value = 1999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1999
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 171

$ws.Range("H123").Value = 28000.273
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 28000.273
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 28000.273
$ws.Range("N123").Value = -32900.273

$ws.Range("H126").Value = 5342.6665
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 5514
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 16542
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -21482

$ws.Range("H132").Value = 3048.5386
$ws.Range("I132").Value = 2763.1
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 8289.299999999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -5759.299999999999
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3255.5293
$ws.Range("I46").Value = 2633.3333
$ws.Range("J46").Value = 3594.9092
$ws.Range("K46").Value = 2633.3333
$ws.Range("L46").Value = 3594.9092
$ws.Range("M46").Value = -2445.3333
$ws.Range("N46").Value = -3970.9092

$ws.Range("H55").Value = 750
$ws.Range("I55").Value = 750
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 750
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -577

$ws.Range("H82").Value = 1432.8
$ws.Range("I82").Value = 1291
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 1291
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -930
$ws.Range("N82").Value = -2722

$ws.Range("H85").Value = 1432.8
$ws.Range("I85").Value = 1291
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 1291
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -43
$ws.Range("N85").Value = -4496

$ws.Range("H132").Value = 4085.6667
$ws.Range("I132").Value = 3559.2222
$ws.Range("J132").Value = 5665
$ws.Range("K132").Value = 10677.6666
$ws.Range("L132").Value = 16995
$ws.Range("M132").Value = -8147.6666
$ws.Range("N132").Value = -22055

$ws.Range("H136").Value = 3191.5
$ws.Range("I136").Value = 2699.889
$ws.Range("J136").Value = 4666.3335
$ws.Range("K136").Value = 8099.667
$ws.Range("L136").Value = 13999.0005
$ws.Range("M136").Value = -5549.667
$ws.Range("N136").Value = -19099.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8396.549999999999
$ws.Range("I81").Value = 1770
$ws.Range("J81").Value = 13818.272
$ws.Range("K81").Value = 3540
$ws.Range("L81").Value = 27636.544
$ws.Range("M81").Value = -2479
$ws.Range("N81").Value = -29758.544

$ws.Range("H84").Value = 8396.549999999999
$ws.Range("I84").Value = 1770
$ws.Range("J84").Value = 13818.272
$ws.Range("K84").Value = 17700
$ws.Range("L84").Value = 138182.72
$ws.Range("M84").Value = -12396
$ws.Range("N84").Value = -148790.72

$ws.Range("H94").Value = 57497.25
$ws.Range("I94").Value = 57497.25
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 57497.25
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -56596.25

$ws.Range("H132").Value = 40249.543
$ws.Range("I132").Value = 52833.223
$ws.Range("J132").Value = 2498.5
$ws.Range("K132").Value = 158499.669
$ws.Range("L132").Value = 7495.5
$ws.Range("M132").Value = -155969.669
$ws.Range("N132").Value = -12555.5
